# PHC.xlsx - BME_DI_SPH sheet: the "partial_oh" helper columns (B and G)
# used to subtract all four expense buckets (labour/contracts/parts and a
# 4th bucket) from the totals. They should only subtract the first bucket
# (labour), i.e. actual_partial_oh = actual_total_exp - actual_labour_exp
# and budgeted_partial_oh = budgeted_total_exp - budgeted_labour_exp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: actual_partial_oh = actual_total_exp (C) - actual_labour_exp (D)
# Column G: budgeted_partial_oh = budgeted_total_exp (H) - budgeted_labour_exp (I)
# Each cell is written individually (rather than one multi-cell Range
# assignment) so every row keeps its own independent formula.
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 2).Formula = "=C$row-D$row"
    $ws.Cells.Item($row, 7).Formula = "=H$row-I$row"
}

# Restore the saved selection state (cell D10 was left selected/active).
$ws.Range("D10").Select()
